$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Fitness 9"
$ws.Range("A3").Value = "Pledge Fitness"
$ws.Range("A4").Value = "Prime Fitness World"
$ws.Range("A5").Value = "Fitness Reloaded"
$ws.Range("A6").Value = "Fitmax Gym"
$ws.Range("A7").Value = "H2O Fitness Pro"
$ws.Range("A8").Value = "Scross Fit Gym"
$ws.Range("A9").Value = "Nawaz Fitness Pro"
$ws.Range("A10").Value = "Brood Fitness - Semi Personal Training Gym"
$ws.Range("A11").Value = "Lidaas Kick Boxing Academy"
